$d = $word.ActiveDocument

$replacements = @(
    @("752×7=", "489×4="),
    @("988×8=", "434×7="),
    @("444×8=", "212×7="),
    @("892×5=", "780×6="),
    @("998×4=", "881×7="),
    @("389×6=", "650×7="),
    @("278×7=", "376×7="),
    @("943×8=", "925×3="),
    @("494×2=", "755×3="),
    @("432×9=", "806×6="),
    @("260×2=", "974×6="),
    @("788×2=", "626×6="),
    @("466×6=", "261×6="),
    @("643×8=", "457×7="),
    @("369×4=", "541×6="),
    @("866×6=", "975×8="),
    @("895×8=", "961×3="),
    @("983×5=", "796×2="),
    @("335×6=", "914×4="),
    @("786×4=", "467×7="),
    @("878×9=", "136×3="),
    @("582×8=", "693×9="),
    @("480×2=", "648×5="),
    @("246×8=", "406×6="),
    @("274×9=", "829×3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
